# mapping details.xlsx - "maaping updated and custom sort corrected"
#
# Adds two new Commodity/Car rows into the existing mapping table:
#   - a new "SBM / SOYBN CK ML SCR" row right under the existing "SBM" row (new row 6)
#   - a new "WCS / COTTONSEED NBXC" row right above the existing "WCS / COTTONSEED,NBXC"
#     row (new row 10)
#
# Existing rows are relocated by copying whole cells (values + formats) bottom-up, which
# keeps the existing style table intact instead of minting throw-away blank-row styles the
# way a plain row/range Insert() would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Move-Row($fromRow, $toRow) {
    $ws.Range("A$fromRow`:B$fromRow").Copy()
    $ws.Range("A$toRow`:B$toRow").PasteSpecial(-4104)  # xlPasteAll
    $excel.CutCopyMode = 0
}

# --- Relocate existing data rows 6..10 down to 7..12 (room for the two new rows) ---
Move-Row 10 12
Move-Row 9 11
Move-Row 8 9
Move-Row 7 8
Move-Row 6 7

# --- Recreate the trailing blank formatted cell at its new location (row 13) ---
$ws.Range("B13").NumberFormat = "@"

# --- Write the two newly inserted rows ---
# (Row 10 is populated before row 6 so the new shared strings land in the same order as
# the author's workbook: "COTTONSEED NBXC" then "SOYBN CK ML SCR".)

# New row 10: duplicate of row 11's formatting (WCS / COTTONSEED,NBXC), new Commodity Name
$ws.Range("A11:B11").Copy()
$ws.Range("A10:B10").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = 0
$ws.Range("A10").Value = "WCS"
$ws.Range("B10").Value = "COTTONSEED NBXC"

# New row 6: duplicate of row 5's formatting (SBM), new Commodity Name
$ws.Range("A5:B5").Copy()
$ws.Range("A6:B6").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = 0
$ws.Range("A6").Value = "SBM"
$ws.Range("B6").Value = "SOYBN CK ML SCR"

# --- Move the selection cursor to match where the author ended up ---
$ws.Range("B17").Select()
